# Loan RBI, Variable Instalments
#
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the "Late" / "heading" / "Outstanding" columns one to the
#   right (N -> O -> P -> Q) to make room for a new "Variable Instalments"
#   style column.
# - Make "Repayment schedule" the active sheet/tab (it was "NewLoanInput"
#   before).
# - Update the selection on "Repayment schedule" to cell S5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Switch to the Repayment schedule sheet (becomes the active/selected tab).
$ws.Activate()

# Insert a new blank column before column N (14th column).
$ws.Columns.Item(14).Insert()

# Leave the new selection where the author left it.
$ws.Range("S5").Select() | Out-Null
